# Entity Balance workbook update:
# Remove the "Morgan" NPC capability row content, rebase Player HP denominator,
# and add a new "Attack"-based NPC table (columns H:J) mirroring the existing
# HP-based table (columns B:D), including a new "Lapahn Jump" entity.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Left-hand table: Player HP value used for ROUNDDOWN(C/$C$2) ---
$ws.Range("C2").Value = 8

# --- New right-hand table header values ---
$ws.Range("H2").Value = "Player HP :"
$ws.Range("I2").Value = 40

$ws.Range("H3").Value = "Entity Name"
$ws.Range("I3").Value = "Attack"
$ws.Range("J3").Value = "# of Hits"

$ws.Range("H4").Value = "Marines"

$ws.Range("H5").Value = "Marine"
$ws.Range("I5").Value = 2
$ws.Range("J5").Formula = '=ROUNDDOWN($I$2/I5, 0)'

$ws.Range("H6").Value = "Marine Captain"
$ws.Range("I6").Value = 6
$ws.Range("J6").Formula = '=ROUNDDOWN($I$2/I6, 0)'

# Row 7 used to mirror "Morgan" on the left table; on the right table it stays blank
# (the NPC capability for Morgan was removed).

$ws.Range("H8").Value = "Pirates"

$ws.Range("H9").Value = "Pirates"
$ws.Range("I9").Value = 2
$ws.Range("J9").Formula = '=ROUNDDOWN($I$2/I9, 0)'

$ws.Range("H10").Value = "Pirate Captain"
$ws.Range("I10").Value = 6
$ws.Range("J10").Formula = '=ROUNDDOWN($I$2/I10, 0)'

$ws.Range("H11").Value = "Fat Pirate"
$ws.Range("I11").Value = 6
$ws.Range("J11").Formula = '=ROUNDDOWN($I$2/I11, 0)'

$ws.Range("H13").Value = "Animals"

$ws.Range("H14").Value = "Kung Fu Dugong"
$ws.Range("I14").Value = 8
$ws.Range("J14").Formula = '=ROUNDDOWN($I$2/I14, 0)'

$ws.Range("H15").Value = "Lapahn"
$ws.Range("I15").Value = 6
$ws.Range("J15").Formula = '=ROUNDDOWN($I$2/I15, 0)'

$ws.Range("H16").Value = "Lapahn Jump"
$ws.Range("I16").Value = 6
$ws.Range("J16").Formula = '=ROUNDDOWN($I$2/I16, 0)'

# --- Mirror the formatting of the left-hand table (B:D) onto the new right-hand
#     table (H:J), row by row (section headers on the right table sit one row
#     higher than their left-table counterparts, since "Morgan" has no analogue). ---
$ws.Range("B3:D3").Copy()
$ws.Range("H3").PasteSpecial(-4122)

$ws.Range("B4:D4").Copy()
$ws.Range("H4").PasteSpecial(-4122)

$ws.Range("B5:D7").Copy()
$ws.Range("H5").PasteSpecial(-4122)

$ws.Range("B9:D9").Copy()
$ws.Range("H8").PasteSpecial(-4122)

$ws.Range("B10:D13").Copy()
$ws.Range("H9").PasteSpecial(-4122)

$ws.Range("B14:D14").Copy()
$ws.Range("H13").PasteSpecial(-4122)

$ws.Range("B15:D16").Copy()
$ws.Range("H14").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# "Lapahn Jump" is a sub-ability, indented under "Lapahn" rather than styled
# like a top-level entity row.
$ws.Range("H16").HorizontalAlignment = -4131
$ws.Range("H16").VerticalAlignment = -4108
$ws.Range("H16").IndentLevel = 2

# --- Merge the new section header cells, matching the left-hand table ---
$ws.Range("H4:J4").Merge()
$ws.Range("H8:J8").Merge()
$ws.Range("H13:J13").Merge()

# --- Remove the stray H27:J27 cells that no longer have any data ---
$ws.Range("H27:J27").ClearContents()

# --- New column widths for the added table / nearby columns ---
$ws.Columns.Item(8).ColumnWidth = 15.736979166666666
$ws.Columns.Item(9).ColumnWidth = 9.166666666666666
$ws.Columns.Item(10).ColumnWidth = 10.736979166666666
$ws.Columns.Item(14).ColumnWidth = 14.877604166666666
$ws.Columns.Item(17).ColumnWidth = 17.166666666666668

# --- Selection moved to C10 in the saved file ---
$ws.Range("C10").Select()
